$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.404.75"
$ws.Range("E2").Value = "  -7.61%  "

# Row 3
$ws.Range("D3").Value = "3.400.14"
$ws.Range("E3").Value = "  -5.28%  "

# Row 4
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.57%  "

# Row 5
$ws.Range("D5").Value = "381.03"
$ws.Range("E5").Value = "  -8.40%  "

# Row 6
$ws.Range("D6").Value = "120.37"
$ws.Range("E6").Value = "  -6.94%  "

# Row 7
$ws.Range("D7").Value = "3.526.03"
$ws.Range("E7").Value = "  -1.57%  "

# Row 8
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  -12.18%  "

# Row 9
$ws.Range("D9").Value = "0.998"
$ws.Range("E9").Value = "  -0.14%  "

# Row 10
$ws.Range("D10").Value = "0.647"
$ws.Range("E10").Value = "  -16.12%  "

# Row 11
$ws.Range("D11").Value = "0.136"
$ws.Range("E11").Value = "  -24.34%  "

# Row 12
$ws.Range("D12").Value = "0.0000280"
$ws.Range("E12").Value = "  -16.47%  "

# Row 13
$ws.Range("D13").Value = "37.68"
$ws.Range("E13").Value = "  -10.87%  "

# Row 14
$ws.Range("D14").Value = "3.905.27"
$ws.Range("E14").Value = "  -6.07%  "

# Row 15
$ws.Range("D15").Value = "8.92"
$ws.Range("E15").Value = "  -9.91%  "

# Row 16
$ws.Range("D16").Value = "0.136"
$ws.Range("E16").Value = "  -3.29%  "

# Row 17
$ws.Range("D17").Value = "3.386.06"
$ws.Range("E17").Value = "  -5.17%  "

# Row 18
$ws.Range("D18").Value = "18.09"
$ws.Range("E18").Value = "  -11.28%  "

# Row 19
$ws.Range("D19").Value = "11.97"
$ws.Range("E19").Value = "  -2.33%  "

# Row 20
$ws.Range("D20").Value = "62.087.41"
$ws.Range("E20").Value = "  -8.05%  "

# Row 21
$ws.Range("D21").Value = "0.986"
$ws.Range("E21").Value = "  -13.33%  "

# Row 22
$ws.Range("D22").Value = "380.75"
$ws.Range("E22").Value = "  -15.77%  "

# Row 23
$ws.Range("D23").Value = "13.11"
$ws.Range("E23").Value = "  -0.20%  "

# Row 24
$ws.Range("D24").Value = "78.49"
$ws.Range("E24").Value = "  -11.81%  "

# Row 25
$ws.Range("D25").Value = "2.72"
$ws.Range("E25").Value = "  -13.40%  "

# Row 26
$ws.Range("D26").Value = "5.17"
$ws.Range("E26").Value = "  +5.70%  "

# Row 27
$ws.Range("D27").Value = "32.04"
$ws.Range("E27").Value = "  -8.45%  "

# Row 28
$ws.Range("D28").Value = "2.89"
$ws.Range("E28").Value = "  -13.66%  "

# Row 29
$ws.Range("D29").Value = "8.52"
$ws.Range("E29").Value = "  -14.69%  "

# Row 30
$ws.Range("D30").Value = "2.56"
$ws.Range("E30").Value = "  -7.74%  "

# Row 31
$ws.Range("D31").Value = "11.45"
$ws.Range("E31").Value = "  -7.24%  "

# Row 32
$ws.Range("D32").Value = "0.106"
$ws.Range("E32").Value = "  -9.34%  "

# Row 33
$ws.Range("D33").Value = "6.50"
$ws.Range("E33").Value = "  -12.13%  "

# Row 34
$ws.Range("E34").Value = "  -0.04%  "

# Row 35
$ws.Range("D35").Value = "0.144"
$ws.Range("E35").Value = "  -10.61%  "

# Row 36
$ws.Range("D36").Value = "53.38"
$ws.Range("E36").Value = "  -5.88%  "

# Row 37
$ws.Range("D37").Value = "35.33"
$ws.Range("E37").Value = "  -13.03%  "

# Row 38
$ws.Range("D38").Value = "0.988"
$ws.Range("E38").Value = "  -1.14%  "

# Row 39
$ws.Range("E39").Value = "  -15.05%  "

# Row 40
$ws.Range("D40").Value = "0.128"
$ws.Range("E40").Value = "  -12.59%  "

# Row 41
$ws.Range("D41").Value = "2.55"
$ws.Range("E41").Value = "  +9.03%  "

# Row 42
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "25.41"
$ws.Range("E42").Value = "  +19.08%  "

# Row 43
$ws.Range("D43").Value = "135.36"
$ws.Range("E43").Value = "  -9.36%  "

# Row 44
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "2.92"
$ws.Range("E44").Value = "  +12.19%  "

# Row 45
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0579"
$ws.Range("E45").Value = "  -23.67%  "

# Row 46
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "2.41"
$ws.Range("E46").Value = "  -12.29%  "

# Row 47
$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").Value = "2.97"
$ws.Range("E47").Value = "  -8.80%  "

# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "1.88"
$ws.Range("E48").Value = "  -4.71%  "

# Row 49
$ws.Range("D49").Value = "3.91"
$ws.Range("E49").Value = "  -9.33%  "

# Row 50
$ws.Range("D50").Value = "2.56"
$ws.Range("E50").Value = "  -16.32%  "

# Row 51
$ws.Range("D51").Value = "0.268"
$ws.Range("E51").Value = "  -14.43%  "
